$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Title paragraph: collapse the three runs
#    "Табела 6.3. ... години (" + "{assessmentYear}" + ") према ... Министарства"
#    into a single run containing just "{header}".
# ---------------------------------------------------------------------------
$titleRange = $d.Paragraphs.Item(1).Range
$titleRange.MoveEnd(1, -1) | Out-Null
$titleRange.Text = "{header}"

# ---------------------------------------------------------------------------
# 2. Table header row: replace the four header-cell captions with the
#    templated placeholders {col1}..{col4}, wrapping the opening brace (and,
#    for col1/col2/col4, the whole placeholder) in DDE-link style bookmarks,
#    mirroring the target markup produced by the authoring tool.
# ---------------------------------------------------------------------------
$tbl = $d.Tables.Item(1)

# --- Column 1: "Редни број" -> "{col1}" ------------------------------------
$cell = $tbl.Cell(1, 1)
$cell.Range.Text = "{col1}"
$cellStart = $tbl.Cell(1, 1).Range.Start
$brace = $d.Range($cellStart, $cellStart + 1)
$d.Bookmarks.Add("__DdeLink__13_2904054838", $brace) | Out-Null

# --- Column 2: "Резултат (назив научног/уметничког резултата)" -> "{col2}" -
$cell = $tbl.Cell(1, 2)
$cell.Range.Text = "{col2}"
$cellStart = $tbl.Cell(1, 2).Range.Start
$whole = $d.Range($cellStart, $cellStart + 6)
$d.Bookmarks.Add("__DdeLink__15_2904054838", $whole) | Out-Null
$brace = $d.Range($cellStart, $cellStart + 1)
$d.Bookmarks.Add("__DdeLink__13_29040548382", $brace) | Out-Null

# --- Column 3: "*Према Правилнику ... (М10, ... M90)" -> "{col3}" ----------
$cell = $tbl.Cell(1, 3)
$cell.Range.Text = "{col3}"
$cellStart = $tbl.Cell(1, 3).Range.Start
$brace = $d.Range($cellStart, $cellStart + 1)
$d.Bookmarks.Add("__DdeLink__13_29040548383", $brace) | Out-Null
$rest = $d.Range($cellStart + 1, $cellStart + 6)
$d.Bookmarks.Add("__DdeLink__17_2904054838", $rest) | Out-Null

# --- Column 4: "Број резултата" -> "{col4}" --------------------------------
$cell = $tbl.Cell(1, 4)
$cell.Range.Text = "{col4}"
$cellStart = $tbl.Cell(1, 4).Range.Start
$whole = $d.Range($cellStart, $cellStart + 6)
$d.Bookmarks.Add("__DdeLink__19_2904054838", $whole) | Out-Null
$brace = $d.Range($cellStart, $cellStart + 1)
$d.Bookmarks.Add("__DdeLink__13_29040548384", $brace) | Out-Null

# ---------------------------------------------------------------------------
# 3. Narrow the last grid column (and its cell) from 2610 dxa to 2609 dxa
#    (130.5pt -> 130.45pt).
# ---------------------------------------------------------------------------
$tbl.Columns.Item(4).Width = 130.45
